# Fix "Recorded By" column (G): swap the order of "System" and the
# recorder's email so that entries read "System, <email>" instead of
# "<email>, System".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1
$lastCol = $usedRange.Column + $usedRange.Columns.Count - 1

# Locate the "Recorded By" column from the header row (row 1) so the
# script is resilient even if columns were reordered.
$recordedByCol = 0
for ($c = 1; $c -le $lastCol; $c++) {
    if ($ws.Cells.Item(1, $c).Value2 -eq "Recorded By") {
        $recordedByCol = $c
        break
    }
}
if ($recordedByCol -eq 0) {
    $recordedByCol = 7
}

$colRange = $ws.Range($ws.Cells.Item(1, $recordedByCol), $ws.Cells.Item($lastRow, $recordedByCol))

$colRange.Replace("dnasr281@gmail.com, System", "System, dnasr281@gmail.com")
